$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview and "Correspond Handoff Datetime"
# on de-de share the same underlying timestamp for the 0721bd8a file, so set
# them to the identical new value.
$wsOverview.Range("G3").Value = "2016-08-29 16:49:52"
$wsDeDe.Range("H3").Value = "2016-08-29 16:49:52"

$wsZhCn.Range("H3").Value = "2016-08-29 16:49:47"
$wsZhCn.Range("K3").Value = "2016-08-29 16:50:21"

$wsDeDe.Range("K3").Value = "2016-08-29 16:50:31"
